# Update column G (K = strikeouts) for rows 2-9 to use true strikeout
# counts instead of the previous "Strike#" values.
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 4
    3 = 4
    4 = 4
    5 = 8
    6 = 7
    7 = 4
    8 = 7
    9 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
